$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted ahead of the existing rows (the
# sheet is sorted most-recent-date-first for this market/category), so
# every existing data row from 499 down to 538 shifts down by one (to
# 500..539) and a brand-new row 499 is populated with the latest reading.
$ws.Rows(499).Insert()

$ws.Range("A499").Value = 5
$ws.Range("B499").Value = "Macroferia Regional de Talca"
$ws.Range("C499").Value = "Maule"
$ws.Range("D499").Value = 45265
$ws.Range("E499").Value = 7
$ws.Range("F499").Value = 100112008
$ws.Range("G499").Value = "Coliflor"
$ws.Range("H499").Value = "Sin especificar"
$ws.Range("I499").Value = "Primera"
$ws.Range("J499").Value = 4000
$ws.Range("K499").Value = 900
$ws.Range("L499").Value = 900
$ws.Range("M499").Value = 900
$ws.Range("N499").Value = "$/unidad"
$ws.Range("O499").Value = "Región del Maule"
$ws.Range("P499").Value = 900
$ws.Range("Q499").Value = 1
$ws.Range("R499").Value = "Hortaliza"
